$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "data/damage_manual/static/hazard/flood100.tif"
$ws.Range("D3").Value = "data/damage_manual/static/hazard/flood100.tif"
$ws.Range("D4").Value = "data/damage_manual/static/hazard/flood1000.tif"
$ws.Range("D5").Value = "data/damage_manual/static/hazard/flood1000.tif"
